$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 1).Value = 0.05426330119371414
$ws.Cells.Item(2, 2).Value = 0.9860075712203979
$ws.Cells.Item(2, 3).Value = 0.02887767739593983
$ws.Cells.Item(2, 4).Value = 0.9953567981719971
$ws.Cells.Item(3, 1).Value = 0.009407783858478069
$ws.Cells.Item(3, 2).Value = 0.9983304738998413
$ws.Cells.Item(3, 3).Value = 0.01511609647423029
$ws.Cells.Item(3, 4).Value = 0.9960899353027344
$ws.Cells.Item(4, 1).Value = 0.004760540090501308
$ws.Cells.Item(4, 2).Value = 0.9987080693244934
$ws.Cells.Item(4, 3).Value = 0.003963137045502663
$ws.Cells.Item(4, 4).Value = 0.9978005886077881
$ws.Cells.Item(5, 1).Value = 0.00272815371863544
$ws.Cells.Item(5, 2).Value = 0.999403715133667
$ws.Cells.Item(5, 3).Value = 0.0005846252315677702
$ws.Cells.Item(6, 1).Value = 0.001506247208453715
$ws.Cells.Item(6, 2).Value = 0.9995229840278625
$ws.Cells.Item(6, 3).Value = 0.0002919467224273831
$ws.Cells.Item(7, 1).Value = 0.001273780129849911
$ws.Cells.Item(7, 2).Value = 0.9995229840278625
$ws.Cells.Item(7, 3).Value = 0.0003336528607178479
$ws.Cells.Item(7, 4).Value = 0.9997556209564209
$ws.Cells.Item(8, 1).Value = 0.0015819794498384
$ws.Cells.Item(8, 2).Value = 0.9994633793830872
$ws.Cells.Item(8, 3).Value = 0.0004527504206635058
$ws.Cells.Item(8, 4).Value = 0.9995112419128418
$ws.Cells.Item(9, 1).Value = 0.0006882870802655816
$ws.Cells.Item(9, 2).Value = 0.9998409748077393
$ws.Cells.Item(9, 3).Value = [double]"8.262062328867614E-05"
$ws.Cells.Item(10, 1).Value = 0.001153822522610426
$ws.Cells.Item(10, 2).Value = 0.9997018575668335
$ws.Cells.Item(10, 3).Value = 0.0001179013997898437
$ws.Cells.Item(10, 4).Value = 1
$ws.Cells.Item(11, 1).Value = 0.000512516766320914
$ws.Cells.Item(11, 2).Value = 0.9998409748077393
$ws.Cells.Item(11, 3).Value = [double]"1.266226445295615E-05"
$ws.Cells.Item(12, 1).Value = 0.0006529099773615599
$ws.Cells.Item(12, 2).Value = 0.9997813701629639
$ws.Cells.Item(12, 3).Value = 0.0002332705189473927
$ws.Cells.Item(12, 4).Value = 0.9997556209564209
$ws.Cells.Item(13, 1).Value = 0.0005341874784789979
$ws.Cells.Item(13, 2).Value = 0.9998608827590942
$ws.Cells.Item(13, 3).Value = [double]"7.043559435260249E-06"
$ws.Cells.Item(14, 1).Value = 0.0006503761978819966
$ws.Cells.Item(14, 2).Value = 0.9997615218162537
$ws.Cells.Item(14, 3).Value = [double]"5.714269354939461E-05"
$ws.Cells.Item(15, 1).Value = 0.0001791017857613042
$ws.Cells.Item(15, 2).Value = 0.9999602437019348
$ws.Cells.Item(15, 3).Value = [double]"1.037123297464859E-06"
$ws.Cells.Item(16, 1).Value = 0.0005436852807179093
$ws.Cells.Item(16, 2).Value = 0.9998012185096741
$ws.Cells.Item(16, 3).Value = [double]"1.303612748415617E-06"
$ws.Cells.Item(17, 1).Value = 0.000437789160059765
$ws.Cells.Item(17, 2).Value = 0.9998409748077393
$ws.Cells.Item(17, 3).Value = [double]"5.712416282221966E-07"
$ws.Cells.Item(18, 1).Value = 0.0007883654325269163
$ws.Cells.Item(18, 2).Value = 0.9997813701629639
$ws.Cells.Item(18, 3).Value = [double]"1.757780637490214E-06"
$ws.Cells.Item(19, 1).Value = 0.0005108626210130751
$ws.Cells.Item(19, 2).Value = 0.9998012185096741
$ws.Cells.Item(19, 3).Value = [double]"7.616019956913078E-07"
$ws.Cells.Item(20, 1).Value = [double]"9.293340553995222E-05"
$ws.Cells.Item(20, 2).Value = 0.9999801516532898
$ws.Cells.Item(20, 3).Value = [double]"4.075154265592573E-06"
$ws.Cells.Item(21, 1).Value = 0.0003835852257907391
$ws.Cells.Item(21, 2).Value = 0.9999006390571594
$ws.Cells.Item(21, 3).Value = [double]"1.660624207033834E-07"
$ws.Cells.Item(22, 1).Value = [double]"2.202675568696577E-05"
$ws.Cells.Item(22, 2).Value = 1
$ws.Cells.Item(22, 3).Value = [double]"1.919956957863178E-07"
$ws.Cells.Item(23, 1).Value = 0.0006493227556347847
$ws.Cells.Item(23, 2).Value = 0.9999006390571594
$ws.Cells.Item(23, 3).Value = [double]"1.299254677178396E-07"
$ws.Cells.Item(24, 1).Value = [double]"1.635597618587781E-05"
$ws.Cells.Item(24, 2).Value = 1
$ws.Cells.Item(24, 3).Value = [double]"5.394887026000106E-08"
$ws.Cells.Item(25, 1).Value = 0.0009091651299968362
$ws.Cells.Item(25, 2).Value = 0.9998807311058044
$ws.Cells.Item(25, 3).Value = [double]"8.584633178543299E-06"
$ws.Cells.Item(25, 4).Value = 1
$ws.Cells.Item(26, 1).Value = 0.0005335803143680096
$ws.Cells.Item(26, 2).Value = 0.9999403953552246
$ws.Cells.Item(26, 3).Value = [double]"8.121498467517085E-06"
$ws.Cells.Item(27, 1).Value = 0.0006532507250085473
$ws.Cells.Item(27, 2).Value = 0.9998807311058044
$ws.Cells.Item(27, 3).Value = [double]"3.167456497976673E-06"
$ws.Cells.Item(28, 1).Value = [double]"8.932907076086849E-05"
$ws.Cells.Item(28, 2).Value = 0.9999602437019348
$ws.Cells.Item(28, 3).Value = [double]"1.095687821361935E-06"
$ws.Cells.Item(29, 1).Value = 0.0001163829147117212
$ws.Cells.Item(29, 2).Value = 0.9999204874038696
$ws.Cells.Item(29, 3).Value = [double]"9.694232971924066E-08"
$ws.Cells.Item(30, 1).Value = 0.000612867355812341
$ws.Cells.Item(30, 2).Value = 0.9998807311058044
$ws.Cells.Item(30, 3).Value = [double]"1.201821078211651E-06"
$ws.Cells.Item(31, 1).Value = 0.0001545879349578172
$ws.Cells.Item(31, 2).Value = 0.9999403953552246
$ws.Cells.Item(31, 3).Value = [double]"8.528054991074896E-07"
$ws.Cells.Item(32, 1).Value = 0.0001932430604938418
$ws.Cells.Item(32, 2).Value = 0.9999602437019348
$ws.Cells.Item(32, 3).Value = [double]"4.237329847001092E-07"
$ws.Cells.Item(33, 1).Value = 0.0005914644571021199
$ws.Cells.Item(33, 2).Value = 0.9998409748077393
$ws.Cells.Item(33, 3).Value = [double]"2.227418320899233E-07"
$ws.Cells.Item(34, 1).Value = [double]"4.286908006179146E-05"
$ws.Cells.Item(34, 2).Value = 0.9999801516532898
$ws.Cells.Item(34, 3).Value = [double]"3.948123890040733E-07"
$ws.Cells.Item(35, 1).Value = [double]"4.72977917524986E-05"
$ws.Cells.Item(35, 2).Value = 0.9999801516532898
$ws.Cells.Item(35, 3).Value = [double]"3.577325287551503E-06"
$ws.Cells.Item(36, 1).Value = [double]"8.08629920356907E-06"
$ws.Cells.Item(36, 2).Value = 1
$ws.Cells.Item(36, 3).Value = [double]"1.246882703753727E-07"
$ws.Cells.Item(37, 1).Value = 0.0001153373013949022
$ws.Cells.Item(37, 2).Value = 0.9999602437019348
$ws.Cells.Item(37, 3).Value = [double]"3.052922181723261E-08"
$ws.Cells.Item(38, 1).Value = 0.0001846843370003626
$ws.Cells.Item(38, 2).Value = 0.9999602437019348
$ws.Cells.Item(38, 3).Value = [double]"2.112530455633532E-05"
$ws.Cells.Item(39, 1).Value = 0.0003542202175594866
$ws.Cells.Item(39, 2).Value = 0.9999403953552246
$ws.Cells.Item(39, 3).Value = [double]"1.246841208057958E-08"
$ws.Cells.Item(40, 1).Value = 0.0001188235764857382
$ws.Cells.Item(40, 2).Value = 0.9999801516532898
$ws.Cells.Item(40, 3).Value = [double]"8.060071365889598E-08"
$ws.Cells.Item(41, 1).Value = [double]"1.567881736264098E-05"
$ws.Cells.Item(41, 2).Value = 1
$ws.Cells.Item(41, 3).Value = [double]"2.770335072455055E-08"
$ws.Cells.Item(42, 1).Value = 0.0002285463851876557
$ws.Cells.Item(42, 2).Value = 0.9999403953552246
$ws.Cells.Item(42, 3).Value = [double]"3.565592976428889E-08"
$ws.Cells.Item(43, 1).Value = [double]"1.176493333332473E-05"
$ws.Cells.Item(43, 2).Value = 1
$ws.Cells.Item(43, 3).Value = [double]"1.04001651735075E-08"
$ws.Cells.Item(44, 1).Value = 0.0005205409834161401
$ws.Cells.Item(44, 2).Value = 0.9999006390571594
$ws.Cells.Item(44, 3).Value = [double]"6.624172499414271E-08"
$ws.Cells.Item(45, 1).Value = [double]"9.480329936195631E-06"
$ws.Cells.Item(45, 2).Value = 1
$ws.Cells.Item(45, 3).Value = [double]"5.681357606590609E-07"
$ws.Cells.Item(46, 1).Value = [double]"8.399745274800807E-05"
$ws.Cells.Item(46, 2).Value = 0.9999801516532898
$ws.Cells.Item(46, 3).Value = [double]"3.10118707602669E-06"
$ws.Cells.Item(47, 1).Value = [double]"6.582752575923223E-06"
$ws.Cells.Item(47, 2).Value = 1
$ws.Cells.Item(47, 3).Value = [double]"1.958312896022107E-05"
$ws.Cells.Item(48, 1).Value = 0.0001552749308757484
$ws.Cells.Item(48, 2).Value = 0.9999801516532898
$ws.Cells.Item(48, 3).Value = [double]"7.133693088690052E-06"
$ws.Cells.Item(49, 1).Value = 0.0004201754927635193
$ws.Cells.Item(49, 2).Value = 0.9999602437019348
$ws.Cells.Item(49, 3).Value = [double]"2.272312560336331E-09"
$ws.Cells.Item(50, 1).Value = [double]"4.054515557072591E-06"
$ws.Cells.Item(50, 3).Value = [double]"1.689669293369889E-09"
$ws.Cells.Item(51, 1).Value = 0.0002574236714281142
$ws.Cells.Item(51, 2).Value = 0.9999403953552246
$ws.Cells.Item(51, 3).Value = [double]"1.187149905490514E-06"
